$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.428.63"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "2.313.28"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'319.14"
$ws.Range("D6").Value = "'103.51"
$ws.Range("E6").Value = "  -2.29%  "
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "'0.611"
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("D10").Value = "'39.76"
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("D11").Value = "'0.0910"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").Value = "'0.107"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").Value = "'0.971"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").Value = "'15.37"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("D16").Value = "2.663.78"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "2.311.66"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").Value = "42.625.05"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").Value = "'73.37"
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("B22").Value = "PancakeSwap"
$ws.Range("C22").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D22").Value = "'3.61"
$ws.Range("E22").Value = "  +4.18%  "
$ws.Range("D23").Value = "'280.71"
$ws.Range("E23").Value = "  +7.43%  "
$ws.Range("D24").Value = "'10.78"
$ws.Range("E24").Value = "  +16.93%  "
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("E28").Value = "  +4.55%  "
$ws.Range("D29").Value = "'22.96"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'35.96"
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("D31").Value = "'165.27"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").Value = "'0.0877"
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("D33").Value = "'5.92"
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("E34").Value = "  +5.61%  "
$ws.Range("D35").Value = "'2.62"
$ws.Range("E35").Value = "  -9.71%  "
$ws.Range("D36").Value = "'0.117"
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("D37").Value = "'4.63"
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("D38").Value = "'0.0363"
$ws.Range("E38").Value = "  +3.86%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "'2.79"
$ws.Range("E39").Value = "  +5.94%  "
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").Value = "'3.70"
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("D41").Value = "'1.50"
$ws.Range("E41").Value = "  +2.50%  "
$ws.Range("D42").Value = "'99.46"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("D43").Value = "'69.67"
$ws.Range("E43").Value = "  -0.32%  "
$ws.Range("D44").Value = "'0.226"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "'12.13"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("D47").Value = "'113.40"
$ws.Range("E47").Value = "  +2.11%  "
$ws.Range("D48").Value = "'79.07"
$ws.Range("E48").Value = "  +7.82%  "
$ws.Range("D49").Value = "'8.98"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").Value = "'5.33"
$ws.Range("E50").Value = "  -1.01%  "
$ws.Range("D51").Value = "1.616.56"
$ws.Range("E51").Value = "  +5.20%  "

Write-Host "Applied 95 cell updates"
